$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price/volume refresh snapshot.
# D (Price) cells can look like plain numbers (e.g. "239.21"), and Excel
# auto-converts such literals to the Number type on assignment. The source
# data keeps these as plain text cells (some prices even use "."  as a
# thousands separator, e.g. "29.338.68"), so force Text format before the
# write and then drop the format override with ClearFormats so no stray
# style index is left behind on the cell.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.338.68"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.87%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.839.89"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -0.67%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.21"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.75%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6257"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.76%  "

$ws.Range("E7").Value = "  +0.07%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07378"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -1.48%  "

$ws.Range("E9").Value = "  -0.97%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.75"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.39%  "

$ws.Range("E11").Value = "  -0.39%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.823.39"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -1.43%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.954"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -1.79%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6638"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -2.77%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001044"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.59%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "81.43"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -1.56%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.243"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.71%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "29.294.68"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.99%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "234.55"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.83%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.24"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.27%  "

$ws.Range("E21").Value = "  +0.07%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.292"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -3.76%  "

$ws.Range("E23").Value = "  +0.06%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "157.31"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.30%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.446"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.97%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1334"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -2.55%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.26"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.00%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.07121"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +6.75%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.479"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +2.64%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.482"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.23%  "

$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.039"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.77%  "

$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.021"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -2.17%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.155"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.67%  "

$ws.Range("E34").Value = "  -1.70%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7016"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.18%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.584"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.78%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01830"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -2.23%  "

$ws.Range("B38").Value = "FraxShare"
$ws.Range("C38").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.816"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.01%  "

$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.784"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -2.12%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.233.65"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -2.53%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9439"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.41%  "

$ws.Range("E42").Value = "  -0.03%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.987.85"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.31%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.23"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.19%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "65.06"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -2.04%  "

$ws.Range("E46").Value = "  -2.65%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.940"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.46%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.683"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -2.96%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.878"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.34%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.1129"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -3.39%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3872"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -2.19%  "
